$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "model_evaluation.ipynb file" -> split into
#   proofErr(spellStart) + "model_evaluation.ipynb" + proofErr(spellEnd) + " file"
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("model_evaluation.ipynb file", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pPrQ = '<w:pPr><w:widowControl w:val="0"/><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr>'
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $pPrQ + '<w:r><w:t xml:space="preserve">What was the response of the model to your domain-specific input in the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>model_evaluation.ipynb</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> file</w:t></w:r><w:r><w:t xml:space="preserve">? </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para1 = $r.Paragraphs(1).Range
$para1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Edit 2: "> relational databases are not well-suited to the needs of the IoT."
#          -> "> data warehousing and data marts are not suited to real-time analytics."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    ">  relational databases are not well-suited to the needs of the IoT.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ">  data warehousing and data marts are not suited to real-time analytics.", 2) | Out-Null

# ---------------------------------------------------------------------------
# Edit 3: "The most important thing to know about IoT ... various tasks."
#          -> "In this presentation, we will discuss ... benefits of real-"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "The most important thing to know about IoT is that it is not a single technology. IoT is a network of devices that are constantly communicating with each other to perform various tasks.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In this presentation, we will discuss the characteristics of real-time analytics, the challenges of real-time analytics and how to overcome these challenges. We will also discuss the benefits of real-", 2) | Out-Null

# ---------------------------------------------------------------------------
# Edit 4: delete the whole paragraph "This network of devices is constantly communic"
# ---------------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("This network of devices is constantly communic", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para4 = $r4.Paragraphs(1).Range
$para4.Delete()

# ---------------------------------------------------------------------------
# Edit 5: "model_finetuning.ipynb file" -> split into
#   proofErr(spellStart) + "model_finetuning.ipynb" + proofErr(spellEnd) + " file"
# ---------------------------------------------------------------------------
$r5 = $d.Content
$r5.Find.Execute("model_finetuning.ipynb file", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pPrQ5 = '<w:pPr><w:widowControl w:val="0"/><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr>'
$xml5 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $pPrQ5 + '<w:r><w:t xml:space="preserve">After fine-tuning the model, what was the response of the model to your domain-specific input in the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>model_finetuning.ipynb</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> file</w:t></w:r><w:r><w:t>?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para5 = $r5.Paragraphs(1).Range
$para5.InsertXML($xml5)

# ---------------------------------------------------------------------------
# Edit 6: "Error" paragraph -> 5 new paragraphs of generated-text transcript
# ---------------------------------------------------------------------------
$r6 = $d.Content
$r6.Find.Execute("Error", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para6 = $r6.Paragraphs(1).Range
$pPr = '<w:pPr><w:widowControl w:val="0"/><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>'
$pPrNoLang = '<w:pPr><w:widowControl w:val="0"/><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>'
$rPrLang = '<w:rPr><w:lang w:val="en-US"/></w:rPr>'

$p1 = '<w:p>' + $pPr + '<w:r>' + $rPrLang + '<w:t>Traditional approaches to data management such as</w:t></w:r></w:p>'
$p2 = '<w:p>' + $pPr + '<w:r>' + $rPrLang + "<w:t>&gt; [{'</w:t></w:r>" + '<w:proofErr w:type="spellStart"/><w:r>' + $rPrLang + '<w:t>generated_text</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r>' + $rPrLang + '<w:t xml:space="preserve">' + "'" + ': ' + "'" + ' relational databases, file systems, and other forms of data storage have been around for decades. However, these approaches are not always suitable for the complex, dynamic, and distributed nature of modern data. In this article, we will explore the benefits </w:t></w:r><w:r>' + $rPrLang + '<w:lastRenderedPageBreak/><w:t>' + "of using a modern data storage solution, such as a distributed file system'}]" + '</w:t></w:r></w:p>'
$p3 = '<w:p>' + $pPr + '</w:p>'
$p4 = '<w:p>' + $pPr + '<w:r>' + $rPrLang + '<w:t>==================================</w:t></w:r></w:p>'
$p5 = '<w:p>' + $pPrNoLang + '</w:p>'

$xml6 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $p1 + $p2 + $p3 + $p4 + $p5 + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para6.InsertXML($xml6)

Write-Output "All edits applied"
